# US Foods_Bakery order sheet: append the new Chocolate Chips line item as
# row 28 (SKU, Item, Quantity, Cost Per, Total Cost), matching the existing
# sheet's convention of storing every value (including numeric-looking
# ones) as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A28:E28")

# Force text storage so "92210", "4", "63.99", "255.96" are not reinterpreted
# as numbers, then restore the default "Normal" style so no stray
# NumberFormat/style is left behind on the cells.
$newRow.NumberFormat = "@"

$ws.Range("A28").Value = "92210"
$ws.Range("B28").Value = "Chocolate Chips 4M - Semi-Sweet"
$ws.Range("C28").Value = "4"
$ws.Range("D28").Value = "63.99"
$ws.Range("E28").Value = "255.96"

$newRow.Style = "Normal"
